# Auto-generated files on 2025-08-27
# Update the hot-stock ranking table (columns A-C, rows 2-21) with the
# refreshed values from the source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "岩山科技"

$ws.Range("B3").Value = "长川科技"
$ws.Range("C3").Value = "华胜天成"

$ws.Range("B4").Value = "麦格米特"
$ws.Range("C4").Value = "寒武纪"

$ws.Range("A5").Value = "北方稀土"
$ws.Range("B5").Value = "英维克"

$ws.Range("A6").Value = "岩山科技"
$ws.Range("B6").Value = "紫光国微"

$ws.Range("A7").Value = "领益智造"
$ws.Range("B7").Value = "科大智能"
$ws.Range("C7").Value = "指南针"

$ws.Range("A8").Value = "剑桥科技"
$ws.Range("B8").Value = "亿纬锂能"
$ws.Range("C8").Value = "吉视传媒"

$ws.Range("B9").Value = "章源钨业"

$ws.Range("A10").Value = "利欧股份"
$ws.Range("B10").Value = "太辰光"
$ws.Range("C10").Value = "领益智造"

$ws.Range("A11").Value = "吉视传媒"
$ws.Range("B11").Value = "隆扬电子"
$ws.Range("C11").Value = "启明信息"

$ws.Range("A12").Value = "新易盛"
$ws.Range("B12").Value = "步步高"
$ws.Range("C12").Value = "剑桥科技"

$ws.Range("A13").Value = "瑞芯微"
$ws.Range("B13").Value = "中油资本"
$ws.Range("C13").Value = "万通发展"

$ws.Range("A14").Value = "华银电力"
$ws.Range("B14").Value = "联化科技"
$ws.Range("C14").Value = "利欧股份"

$ws.Range("A15").Value = "合力泰"
$ws.Range("B15").Value = "同花顺"

$ws.Range("A16").Value = "启明信息"
$ws.Range("B16").Value = "平潭发展"
$ws.Range("C16").Value = "华银电力"

$ws.Range("A17").Value = "英维克"
$ws.Range("B17").Value = "润建股份"
$ws.Range("C17").Value = "新易盛"

$ws.Range("A18").Value = "麦格米特"
$ws.Range("B18").Value = "蔚蓝锂芯"

$ws.Range("B19").Value = "欣旺达"
$ws.Range("C19").Value = "方正科技"

$ws.Range("A20").Value = "步步高"
$ws.Range("B20").Value = "科华数据"
$ws.Range("C20").Value = "国光连锁"

$ws.Range("A21").Value = "博通集成"
$ws.Range("B21").Value = "全志科技"
